$d = $word.ActiveDocument

# --- Fix typo "Deumeurant" -> "Demeurant" and merge the split runs
#     ("Deumeurant au : " + "[ADDR]") into a single run "Demeurant au : [ADDR]".
$find = $d.Content.Find
$find.Execute("Deumeurant au : [ADDR]", $false, $false, $false, $false, $false, `
               $true, 1, $false, "Demeurant au : [ADDR]", 2)

# --- Add the missing "ListLabel 23" character style (mirrors ListLabel22's
#     formatting: complex-script font Symbol, size 12pt/24 half-points).
$newStyle = $d.Styles.Add("ListLabel 23", 2)
$newStyle.Font.NameBi = "Symbol"
$newStyle.Font.Size = 12
$newStyle.Font.SizeBi = 12
$newStyle.QuickStyle = $true
